$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.165.40'
$ws.Range('E2').Value = '  -0.81%  '
$ws.Range('D3').Value = '2.693.42'
$ws.Range('E3').Value = '  +1.58%  '
$ws.Range('E4').Value = '  +0.01%  '
$c = $ws.Range('D5')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '599.39'
$c.Style = $origStyle
$ws.Range('E5').Value = '  -1.77%  '
$c = $ws.Range('D6')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '174.86'
$c.Style = $origStyle
$ws.Range('E6').Value = '  -4.41%  '
$ws.Range('E7').Value = '  +0.01%  '
$c = $ws.Range('D8')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.524'
$c.Style = $origStyle
$ws.Range('E8').Value = '  -1.10%  '
$ws.Range('D9').Value = '2.691.92'
$ws.Range('E9').Value = '  +1.58%  '
$ws.Range('E10').Value = '  -6.71%  '
$ws.Range('E11').Value = '  +1.96%  '
$c = $ws.Range('D12')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.357'
$c.Style = $origStyle
$ws.Range('E12').Value = '  +0.92%  '
$ws.Range('E13').Value = '  -2.45%  '
$ws.Range('D14').Value = '3.180.79'
$ws.Range('E14').Value = '  +1.44%  '
$c = $ws.Range('D15')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.0000185'
$c.Style = $origStyle
$ws.Range('E15').Value = '  -5.59%  '
$ws.Range('D16').Value = '72.026.27'
$ws.Range('E16').Value = '  -0.68%  '
$c = $ws.Range('D17')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '26.27'
$c.Style = $origStyle
$ws.Range('E17').Value = '  -2.58%  '
$ws.Range('D18').Value = '2.686.23'
$ws.Range('E18').Value = '  +1.48%  '
$c = $ws.Range('D19')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '12.25'
$c.Style = $origStyle
$ws.Range('E19').Value = '  +4.79%  '
$c = $ws.Range('D20')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '8.16'
$c.Style = $origStyle
$ws.Range('E20').Value = '  +2.32%  '
$c = $ws.Range('D21')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '372.32'
$c.Style = $origStyle
$ws.Range('E21').Value = '  -3.63%  '
$c = $ws.Range('D22')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '4.18'
$c.Style = $origStyle
$ws.Range('E22').Value = '  -1.12%  '
$c = $ws.Range('D23')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.00'
$c.Style = $origStyle
$ws.Range('E23').Value = '  -1.50%  '
$c = $ws.Range('D24')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '72.35'
$c.Style = $origStyle
$ws.Range('E24').Value = '  -1.71%  '
$ws.Range('E25').Value = '  +0.03%  '
$c = $ws.Range('D26')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '4.35'
$c.Style = $origStyle
$ws.Range('E26').Value = '  -2.79%  '
$c = $ws.Range('D27')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '9.81'
$c.Style = $origStyle
$ws.Range('E27').Value = '  -2.29%  '
$ws.Range('D28').Value = '2.828.41'
$ws.Range('E28').Value = '  +1.49%  '
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('D30').Value = '0.0₃0979'
$ws.Range('E30').Value = '  -0.27%  '
$c = $ws.Range('D31')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '8.07'
$c.Style = $origStyle
$ws.Range('E31').Value = '  -0.89%  '
$c = $ws.Range('D32')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '504.05'
$c.Style = $origStyle
$ws.Range('E32').Value = '  -8.29%  '
$ws.Range('E33').Value = '  -3.78%  '
$ws.Range('E34').Value = '  -1.64%  '
$c = $ws.Range('D36')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '163.91'
$c.Style = $origStyle
$ws.Range('E36').Value = '  -0.19%  '
$c = $ws.Range('D37')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '19.67'
$c.Style = $origStyle
$ws.Range('E37').Value = '  +1.26%  '
$c = $ws.Range('D38')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '19.10'
$c.Style = $origStyle
$ws.Range('E38').Value = '  -0.21%  '
$ws.Range('E39').Value = '  -2.95%  '
$ws.Range('E40').Value = '  -5.22%  '
$ws.Range('E41').Value = '  -5.60%  '
$ws.Range('E42').Value = '  -0.09%  '
$c = $ws.Range('D43')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '5.03'
$c.Style = $origStyle
$ws.Range('E43').Value = '  -2.72%  '
$ws.Range('E44').Value = '  -4.21%  '
$c = $ws.Range('D45')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.333'
$c.Style = $origStyle
$ws.Range('E45').Value = '  -1.17%  '
$c = $ws.Range('D46')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '157.05'
$c.Style = $origStyle
$ws.Range('E46').Value = '  +3.12%  '
$c = $ws.Range('D47')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '39.52'
$c.Style = $origStyle
$ws.Range('E47').Value = '  -0.40%  '
$c = $ws.Range('D48')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.565'
$c.Style = $origStyle
$ws.Range('E48').Value = '  +3.82%  '
$ws.Range('E49').Value = '  +0.64%  '
$c = $ws.Range('D50')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '1.76'
$c.Style = $origStyle
$ws.Range('E50').Value = '  +2.36%  '
$c = $ws.Range('D51')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.0758'
$c.Style = $origStyle
$ws.Range('E51').Value = '  -0.82%  '
